# Insert a new weekly record at row 73, pushing the existing rows 73:137
# down to 74:138, then populate the freshly inserted row with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 73 (shifts 73:137 -> 74:138).
$ws.Rows.Item(73).Insert()

# Populate the new row 73 with the latest weekly price record.
$ws.Range("A73").Value = 1
$ws.Range("B73").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C73").Value = "Arica y Parinacota"
$ws.Range("D73").Value = 44944
$ws.Range("E73").Value = 15
$ws.Range("F73").Value = 100114001
$ws.Range("G73").Value = "Papa"
$ws.Range("H73").Value = "Patagonia"
$ws.Range("I73").Value = "1a (cosecha)"
$ws.Range("J73").Value = 1000
$ws.Range("K73").Value = 14000
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = 14500
$ws.Range("N73").Value = "$/saco 25 kilos"
$ws.Range("O73").Value = "Región del Maule"
$ws.Range("P73").Value = 580
$ws.Range("Q73").Value = 25
$ws.Range("R73").Value = "Hortaliza"
